$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting count/mean/std/... from A:H to B:I
$ws.Columns.Item(1).Insert()

# Copy the header style from the (now shifted) "count" header cell (B1) into the
# new A1 cell so the new header matches the existing bold/bordered header look,
# then set its text to "Feature".
$ws.Range("B1").Copy($ws.Range("A1"))
$ws.Range("A1").Value = "Feature"

# Fill in the feature name for every data row (rows 2-9), in the same order as
# the existing statistics rows.
$ws.Range("A2").Value = "Pregnancies"
$ws.Range("A3").Value = "BloodPressure"
$ws.Range("A4").Value = "SkinThickness"
$ws.Range("A5").Value = "Glucose"
$ws.Range("A6").Value = "Insulin"
$ws.Range("A7").Value = "BMI"
$ws.Range("A8").Value = "DiabetesPedigreeFunction"
$ws.Range("A9").Value = "Age"
